# ===========================================================================
# Edit: visual changes on coordinates and message of length (in Spanish)
#
# - LOCALIZACION: update Norte/Este coordinates on row 3
# - LINEA: populate rows 2-3 (point table) that previously only had headers
# - CIRCULO: add row 3 (second point) to the point table
# - P_DIST_ANG_G:   extend the point table from 37 to 74 points (rows 39-75)
# - P_DIST_ANG_UTM: extend the point table from 37 to 74 points (rows 39-75)
# ===========================================================================

$wb = $excel.ActiveWorkbook

# --- "LOCALIZACION": update Norte (B3) / Este (C3) for the existing 2nd point ---
$wsLoc = $wb.Worksheets.Item("LOCALIZACION")
$wsLoc.Cells.Item(3,2).Value = 1160395.180980101
$wsLoc.Cells.Item(3,3).Value = 754726.4470557353

# --- "LINEA": table was empty (headers only); add the two data rows ---
$wsLin = $wb.Worksheets.Item("LINEA")
# Reuse the index-column style (bold + border) already used on other sheets
$wsLin.Cells.Item(1,2).Copy()
$wsLin.Range("A2:A3").PasteSpecial(-4122)
$lineaData = @(
  @(0,1160227.828983885,731000.3625090888,19),
  @(1,1160395.180980101,754726.4470557353,19)
)
for ($i = 0; $i -lt $lineaData.Length; $i++) {
  $r = 2 + $i
  $row = $lineaData[$i]
  for ($c = 1; $c -le 4; $c++) {
    $wsLin.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}

# --- "CIRCULO": table already had row 2; add row 3 (second point) ---
$wsCir = $wb.Worksheets.Item("CIRCULO")
$wsCir.Cells.Item(2,1).Copy()
$wsCir.Range("A3").PasteSpecial(-4122)
$wsCir.Cells.Item(3,1).Value = 1
$wsCir.Cells.Item(3,2).Value = 1160395.180980101
$wsCir.Cells.Item(3,3).Value = 754726.4470557353
$wsCir.Cells.Item(3,4).Value = 19

# --- "P_DIST_ANG_G": append points 37-73 (rows 39-75), columns A:E ---
$wsG = $wb.Worksheets.Item("P_DIST_ANG_G")
$wsG.Cells.Item(38,1).Copy()
$wsG.Range("A39:A75").PasteSpecial(-4122)
$gData = @(
  @(37,10.53368243087264,-66.67279722222223,0,5),
  @(38,10.52857664048812,-66.66565742994513,10,4.5),
  @(39,10.52253203210929,-66.66029733755495,20,4),
  @(40,10.52377437938858,-66.65223932949182,30,4.5),
  @(41,10.52317278993719,-66.64343203779517,40,5),
  @(42,10.50608875571192,-66.65180074488707,50,3),
  @(43,10.50224051548196,-66.64906066673028,60,3),
  @(44,10.49798288788125,-66.64704190504889,70,3),
  @(45,10.49344525284491,-66.64580577478999,80,3),
  @(46,10.48876549328599,-66.64538980271131,90,3),
  @(47,10.48330543959474,-66.64130823190406,100,3.5),
  @(48,10.4783191109526,-66.64360971725871,110,3.4),
  @(49,10.47079879745071,-66.64115169627692,120,4),
  @(50,10.47202067151893,-66.6525028616766,130,2.9),
  @(51,10.47156259440271,-66.65811707975161,140,2.5),
  @(52,10.46542745785906,-66.65909454505076,150,3),
  @(53,10.46344232221885,-66.66342409903044,160,3),
  @(54,10.45337988958488,-66.66645228229544,170,4),
  @(55,10.45283405530189,-66.67279722222223,180,4),
  @(56,10.45161054743679,-66.67945937114798,190,4.2),
  @(57,10.44824758408166,-66.68779348463246,200,4.8),
  @(58,10.45764759146022,-66.69106699999932,210,4),
  @(59,10.46123982833683,-66.69628466790221,220,4),
  @(60,10.47144320048521,-66.69379134917509,230,3),
  @(61,10.47529105778028,-66.69653171082275,240,3),
  @(62,10.48262136782188,-66.68996658132221,250,2),
  @(63,10.48564634470289,-66.69079106702196,260,2),
  @(64,10.48876653629103,-66.6819330287464,270,1),
  @(65,10.49110624432473,-66.68629284390344,280,1.5),
  @(66,10.49183897075209,-66.68138215743522,290,1),
  @(67,10.49550381123731,-66.68466524144269,300,1.5),
  @(68,10.49454084948603,-66.67979578664722,310,1),
  @(69,10.50941066445419,-66.69041554795545,320,3),
  @(70,10.51988469886784,-66.69107067447369,330,4),
  @(71,10.51409073659551,-66.68217187951728,340,3),
  @(72,10.52415331794662,-66.67914361331242,350,4),
  @(73,10.53368243087264,-66.67279722222223,0,5)
)
for ($i = 0; $i -lt $gData.Length; $i++) {
  $r = 39 + $i
  $row = $gData[$i]
  for ($c = 1; $c -le 5; $c++) {
    $wsG.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}

# --- "P_DIST_ANG_UTM": append points 37-73 (rows 39-75), columns A:F ---
$wsU = $wb.Worksheets.Item("P_DIST_ANG_UTM")
$wsU.Cells.Item(38,1).Copy()
$wsU.Range("A39:A75").PasteSpecial(-4122)
$uData = @(
  @(37,1165365.234029098,754689.5995655446,19,0,5),
  @(38,1164806.078225782,755475.5951735781,19,10,4.5),
  @(39,1164141.596497289,756067.5138806826,19,20,4),
  @(40,1164285.662934859,756948.8521218006,19,30,4.5),
  @(41,1164226.329457347,757913.7696892675,19,40,5),
  @(42,1162329.031729849,757011.5230327214,19,50,3),
  @(43,1161905.451676256,757314.7711895413,19,60,3),
  @(44,1161435.980664305,757539.3731267486,19,70,3),
  @(45,1160934.883913018,757678.5037528875,19,80,3),
  @(46,1160417.387770457,757727.9352189028,19,90,3),
  @(47,1159816.554667152,758179.4570319937,19,100,3.5),
  @(48,1159262.905639278,757931.535880057,19,110,3.4),
  @(49,1158432.761823468,758206.9668368772,19,120,4),
  @(50,1158558.682693018,756962.769526904,19,130,2.9),
  @(51,1158503.418876175,756348.2779553839,19,140,2.5),
  @(52,1157823.748876176,756246.2686062483,19,150,3),
  @(53,1157600.569556579,755773.7162540231,19,160,3),
  @(54,1156484.675957531,755450.2993892395,19,170,4),
  @(55,1156419.146960565,754755.8128866819,19,180,4),
  @(56,1156278.391168782,754027.1357930441,19,190,4.2),
  @(57,1155899.576852499,753117.0675596655,19,200,4.8),
  @(58,1156937.074284995,752750.9246191554,19,210,4),
  @(59,1157330.380457669,752176.5716975678,19,220,4),
  @(60,1158461.389022815,752441.3820568862,19,230,3),
  @(61,1158884.962101866,752138.1522124861,19,240,3),
  @(62,1159701.336216117,752851.1875440651,19,250,2),
  @(63,1160035.391197848,752758.4391027458,19,260,2),
  @(64,1160387.779281084,753725.963648667,19,270,1),
  @(65,1160643.156703452,753246.6078076609,19,280,1.5),
  @(66,1160728.194916033,753783.7844191547,19,290,1),
  @(67,1161131.067059369,753421.2649247225,19,300,1.5),
  @(68,1161028.444888703,753955.3033956827,19,310,1),
  @(69,1162665.247817952,752780.2405421536,19,320,3),
  @(70,1163823.68705968,752699.9738486945,19,330,4),
  @(71,1163189.760149787,753679.1344290869,19,340,3),
  @(72,1164305.659088898,754002.5005536918,19,350,4),
  @(73,1165365.234029098,754689.5995655446,19,0,5)
)
for ($i = 0; $i -lt $uData.Length; $i++) {
  $r = 39 + $i
  $row = $uData[$i]
  for ($c = 1; $c -le 6; $c++) {
    $wsU.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}

